# Weekly cryptos-list refresh: update Price (D) and Volume(1h) (E) figures
# for every listed coin, and replace a few coins whose ranking shifted
# (rows 30-32 rotate PEPE/PancakeSwap/Aptos; row 51 swaps VeChain -> Hedera).
#
# Price cells are stored as text in the source sheet (e.g. "58.909.50",
# "5.60", "0.0930") - several of them look like plain decimals to Excel's
# auto-typing, which would silently convert them to numbers and drop
# trailing zeros (e.g. "5.60" -> 5.6). To keep them as text exactly as
# authored, each Price cell is force-formatted as Text before the write,
# then ClearFormats() removes the now-unneeded style index again so the
# cell's formatting is left exactly as it was found.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-PriceText "D2" "58.909.50"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3 - Ethereum
Set-PriceText "D3" "2.501.21"
$ws.Range("E3").Value = "  +2.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-PriceText "D5" "540.44"
$ws.Range("E5").Value = "  +1.54%  "

# Row 6 - Solana
Set-PriceText "D6" "143.77"
$ws.Range("E6").Value = "  -2.32%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - XRP
Set-PriceText "D8" "0.571"
$ws.Range("E8").Value = "  +0.63%  "

# Row 9 - LidoStakedEther
Set-PriceText "D9" "2.524.93"
$ws.Range("E9").Value = "  +2.43%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.32%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.49%  "

# Row 12 - Toncoin
Set-PriceText "D12" "5.60"
$ws.Range("E12").Value = "  +5.89%  "

# Row 13 - Cardano
Set-PriceText "D13" "0.355"
$ws.Range("E13").Value = "  +1.46%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-PriceText "D14" "2.945.20"
$ws.Range("E14").Value = "  +2.24%  "

# Row 15 - Avalanche
Set-PriceText "D15" "23.61"
$ws.Range("E15").Value = "  -1.14%  "

# Row 16 - WrappedBTC
Set-PriceText "D16" "58.846.96"
$ws.Range("E16").Value = "  -0.55%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +1.62%  "

# Row 18 - WrappedEther
Set-PriceText "D18" "2.519.79"
$ws.Range("E18").Value = "  +0.63%  "

# Row 19 - Chainlink
Set-PriceText "D19" "11.21"
$ws.Range("E19").Value = "  +1.20%  "

# Row 20 - Polkadot
Set-PriceText "D20" "4.28"
$ws.Range("E20").Value = "  -1.14%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "324.90"
$ws.Range("E21").Value = "  +0.69%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +3.26%  "

# Row 23 - Uniswap
Set-PriceText "D23" "5.79"
$ws.Range("E23").Value = "  +0.96%  "

# Row 24 - Litecoin
Set-PriceText "D24" "61.88"
$ws.Range("E24").Value = "  +2.49%  "

# Row 25 - Polygon
Set-PriceText "D25" "0.439"
$ws.Range("E25").Value = "  -4.69%  "

# Row 26 - Kaspa
Set-PriceText "D26" "0.163"
$ws.Range("E26").Value = "  +1.11%  "

# Row 27 - WrappedeETH
Set-PriceText "D27" "2.618.61"
$ws.Range("E27").Value = "  +2.51%  "

# Row 28 - Binance-PegBSC-USD
Set-PriceText "D28" "0.992"
$ws.Range("E28").Value = "  +1.62%  "

# Row 29 - InternetComputer(DFINITY)
Set-PriceText "D29" "7.78"
$ws.Range("E29").Value = "  +1.16%  "

# Row 30 - was PEPE, now PancakeSwap
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-PriceText "D30" "1.81"
$ws.Range("E30").Value = "  +0.01%  "

# Row 31 - was Aptos, now PEPE
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-PriceText "D31" "0.0₃0773"
$ws.Range("E31").Value = "  +0.62%  "

# Row 32 - was PancakeSwap, now Aptos
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-PriceText "D32" "6.66"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33 - Fetch.AI
Set-PriceText "D33" "1.19"
$ws.Range("E33").Value = "  -4.07%  "

# Row 34 - USDe
Set-PriceText "D34" "0.996"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35 - Monero
$ws.Range("E35").Value = "  +0.96%  "

# Row 36 - ImmutableX
Set-PriceText "D36" "1.44"
$ws.Range("E36").Value = "  +3.94%  "

# Row 37 - EthereumClassic
Set-PriceText "D37" "18.65"
$ws.Range("E37").Value = "  +1.75%  "

# Row 38 - NEARProtocol
Set-PriceText "D38" "4.35"
$ws.Range("E38").Value = "  -3.89%  "

# Row 39 - Stacks
Set-PriceText "D39" "1.59"
$ws.Range("E39").Value = "  -7.59%  "

# Row 40 - RenderToken
Set-PriceText "D40" "5.68"
$ws.Range("E40").Value = "  -2.10%  "

# Row 41 - OKB
Set-PriceText "D41" "36.86"
$ws.Range("E41").Value = "  +0.62%  "

# Row 42 - Bittensor
Set-PriceText "D42" "295.83"
$ws.Range("E42").Value = "  -4.92%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  +0.04%  "

# Row 44 - SuiNetwork
Set-PriceText "D44" "0.821"
$ws.Range("E44").Value = "  -0.90%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  -0.06%  "

# Row 46 - Mantle
Set-PriceText "D46" "0.603"
$ws.Range("E46").Value = "  +3.48%  "

# Row 47 - WhiteBITCoin
$ws.Range("E47").Value = "  +0.53%  "

# Row 48 - Stellar
Set-PriceText "D48" "0.0930"
$ws.Range("E48").Value = "  -0.68%  "

# Row 49 - Aave
Set-PriceText "D49" "123.88"
$ws.Range("E49").Value = "  +1.83%  "

# Row 50 - EnergySwap
Set-PriceText "D50" "18.59"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51 - was VeChain, now Hedera
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-PriceText "D51" "0.0514"
$ws.Range("E51").Value = "  -1.66%  "
